# Apply the PtX demand update: insert "Fossil Gases" and "Fossil Liquids" category rows
# into each year block (2030, 2040, 2050) and refresh all labels/values accordingly.
# (Commit: "Updated Outputs folder with correct categories annd Aviation")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert the 6 new category rows. Processed top-to-bottom so each target
# row index is already correct (rows below shift down automatically on Insert()).
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(34).Insert()

# Step 2: (re)write every data row (2-37) with its final label/year/values.
$data = @{}
$data[2] = @("Hydrogen", 2030, $null, $null, $null, 0.0006048130745747752, $null, (4.633996662316279 * [Math]::Pow(10, -9)), 0.0002300625028332096, $null, $null)
$data[3] = @("Methanol", 2030, $null, (7.4739721563973 * [Math]::Pow(10, -5)), $null, $null, $null, $null, $null, $null, $null)
$data[4] = @("Ammonia", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[5] = @("Synthetic Gases", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[6] = @("Biogenic Gases", 2030, $null, $null, 0.0003621210925501598, 0.0001346379081153953, $null, $null, (3.700243579827917 * [Math]::Pow(10, -5)), $null, $null)
$data[7] = @("Fossil Gases", 2030, $null, $null, $null, 0.002327188381595282, $null, $null, 0.0001638215923338511, $null, $null)
$data[8] = @("Synthetic Liquids", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[9] = @("Biogenic Liquids", 2030, $null, $null, $null, 0.0104205708008055, (6.059328194027003 * [Math]::Pow(10, -5)), 0.0089490738962905, 0.007806621961733901, (4.735454024036412 * [Math]::Pow(10, -5)), 0.005130419007949293)
$data[10] = @("Fossil Liquids", 2030, $null, $null, $null, 0.1088088360683084, 0.0004267407160212, 0.0811276643944179, 0.0488992894614437, 0.0002866371317003, 0.0501792817045569)
$data[11] = @("Biomass [Solid]", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[12] = @("Renewable Energy Carrier", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[13] = @("Overall Demand", 2030, $null, (7.4739721563973 * [Math]::Pow(10, -5)), 0.0003621210925501598, 0.1222960462333993, 0.00048733399796147, 0.09007674292470506, 0.05713679795414294, 0.0003339916719406641, 0.05530970071250619)
$data[14] = @("Hydrogen", 2040, $null, $null, $null, 0.002937579253255237, $null, (3.879167034891512 * [Math]::Pow(10, -7)), 0.0003272092465510165, $null, $null)
$data[15] = @("Methanol", 2040, $null, (6.72894938918624 * [Math]::Pow(10, -5)), $null, $null, $null, $null, $null, $null, $null)
$data[16] = @("Ammonia", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[17] = @("Synthetic Gases", 2040, $null, $null, $null, (1.022983814688694 * [Math]::Pow(10, -9)), $null, $null, (9.85606236653328 * [Math]::Pow(10, -11)), $null, $null)
$data[18] = @("Biogenic Gases", 2040, $null, $null, 0.001482638685480181, 0.0001789569760776939, $null, $null, (6.776706688190254 * [Math]::Pow(10, -5)), $null, $null)
$data[19] = @("Fossil Gases", 2040, $null, $null, $null, 0.001218321656675916, $null, $null, 0.0001727691517342355, $null, $null)
$data[20] = @("Synthetic Liquids", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[21] = @("Biogenic Liquids", 2040, $null, $null, $null, 0.0042125048287027, (9.877704592898886 * [Math]::Pow(10, -5)), 0.0111129260681684, 0.0051060944550518, (5.752808082768649 * [Math]::Pow(10, -5)), 0.005809731408513409)
$data[22] = @("Fossil Liquids", 2040, $null, $null, $null, 0.0287835849327677, 0.0004587699952913, 0.07658536659775869, 0.021623831135409, 0.0002543648103777, 0.048663694730089)
$data[23] = @("Biomass [Solid]", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[24] = @("Renewable Energy Carrier", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[25] = @("Overall Demand", 2040, $null, (6.72894938918624 * [Math]::Pow(10, -5)), 0.001482638685480181, 0.03733094867046306, 0.0005575470412202888, 0.08769868058263058, 0.02729767115418858, 0.0003118928912053865, 0.05447342613860241)
$data[26] = @("Hydrogen", 2050, $null, $null, $null, 0.0040596457540099, $null, (6.574875663464897 * [Math]::Pow(10, -7)), 0.0005199128098262, $null, $null)
$data[27] = @("Methanol", 2050, $null, (5.885281635630735 * [Math]::Pow(10, -5)), $null, $null, $null, $null, $null, $null, $null)
$data[28] = @("Ammonia", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[29] = @("Synthetic Gases", 2050, $null, $null, $null, (5.047263276392111 * [Math]::Pow(10, -9)), $null, $null, (2.230309643506691 * [Math]::Pow(10, -9)), $null, $null)
$data[30] = @("Biogenic Gases", 2050, $null, $null, 0.003799263951621619, (2.458547981493456 * [Math]::Pow(10, -5)), $null, $null, (1.976394254763757 * [Math]::Pow(10, -5)), $null, $null)
$data[31] = @("Fossil Gases", 2050, $null, $null, $null, (5.170136468685901 * [Math]::Pow(10, -5)), $null, $null, (6.60249994654634 * [Math]::Pow(10, -5)), $null, $null)
$data[32] = @("Synthetic Liquids", 2050, $null, $null, $null, (3.635139314769846 * [Math]::Pow(10, -11)), (3.637424715607438 * [Math]::Pow(10, -12)), (4.102616117173518 * [Math]::Pow(10, -10)), (8.676915898679511 * [Math]::Pow(10, -11)), (3.93999798732349 * [Math]::Pow(10, -13)), (4.011756553728407 * [Math]::Pow(10, -10)))
$data[33] = @("Biogenic Liquids", 2050, $null, $null, $null, 0.0004114165352042191, 0.0001764972386482, 0.0148291746122799, 0.0013122610033398, (7.399022180948297 * [Math]::Pow(10, -5)), 0.008274360287381074)
$data[34] = @("Fossil Liquids", 2050, $null, $null, $null, 0.002196610500836, 0.0004137911066812, 0.0684978654615489, 0.0038819883370274, 0.0002189447388395, 0.0453483827234034)
$data[35] = @("Biomass [Solid]", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[36] = @("Renewable Energy Carrier", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[37] = @("Overall Demand", 2050, $null, (5.885281635630735 * [Math]::Pow(10, -5)), 0.003799263951621619, 0.006743964718166583, 0.0005902883489668247, 0.08332769797165675, 0.005799953409285303, 0.0002929349610429827, 0.05362274341196013)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le 11; $c++) {
        $v = $rowVals[$c-1]
        if ($v -eq $null) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $v
        }
    }
}

Write-Output "PtX demand table updated: rows 2-37 written."
